$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename personas in column B (row order matters: B3 is edited first so the
# vacated "Yuri" shared-string slot is reused in place for "Ana Valéria",
# while B4/B5 pick up brand new shared-string entries).
$ws.Range("B3").Value = "Ana Valéria"
$ws.Range("B4").Value = "Gustavo Rodrigues"
$ws.Range("B5").Value = "Janaina Rodrigues"

# Column B needs to widen to fit the new (longer) names.
$ws.Columns.Item(2).ColumnWidth = 15.33

# Move the active selection to C8 (matches the post-edit cursor position).
[void]$ws.Range("C8").Select()
